# Refresh ligand/receptor expression metrics for Mcam-Mcam with updated TPM values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 92.64038833333332
$ws.Range("H2").Value = 277.921165
$ws.Range("I2").Value = 0.7451295270557885
$ws.Range("J2").Value = 0.7451295270557885
$ws.Range("M2").Value = 92.64038833333332
$ws.Range("N2").Value = 277.921165
$ws.Range("O2").Value = 0.7451295270557885
$ws.Range("P2").Value = 0.7451295270557885
$ws.Range("Q2").Value = 8582.2415505508
$ws.Range("R2").Value = 77240.17395495721
$ws.Range("S2").Value = 0.555218012090383
$ws.Range("T2").Value = 0.555218012090383
# Row 3
$ws.Range("G3").Value = 92.64038833333332
$ws.Range("H3").Value = 277.921165
$ws.Range("I3").Value = 0.7451295270557885
$ws.Range("J3").Value = 0.7451295270557885
$ws.Range("N3").Value = 6.766394999999999
$ws.Range("O3").Value = 0.01814126213173672
$ws.Range("P3").Value = 0.01814126213173672
$ws.Range("Q3").Value = 208.9471534722416
$ws.Range("R3").Value = 1880.524381250175
$ws.Range("S3").Value = 0.01351759007241607
$ws.Range("T3").Value = 0.01351759007241607
# Row 4
$ws.Range("G4").Value = 92.64038833333332
$ws.Range("H4").Value = 277.921165
$ws.Range("I4").Value = 0.7451295270557885
$ws.Range("J4").Value = 0.7451295270557885
$ws.Range("M4").Value = 28.72545833333334
$ws.Range("N4").Value = 86.17637500000001
$ws.Range("O4").Value = 0.2310459570329316
$ws.Range("P4").Value = 0.2310459570329316
$ws.Range("Q4").Value = 2661.137615052986
$ws.Range("R4").Value = 23950.23853547687
$ws.Range("S4").Value = 0.1721591646921003
$ws.Range("T4").Value = 0.1721591646921004
# Row 5
$ws.Range("G5").Value = 92.64038833333332
$ws.Range("H5").Value = 277.921165
$ws.Range("I5").Value = 0.7451295270557885
$ws.Range("J5").Value = 0.7451295270557885
$ws.Range("M5").Value = 0.706587
$ws.Range("N5").Value = 2.119761
$ws.Range("O5").Value = 0.005683253779543222
$ws.Range("P5").Value = 0.005683253779543222
$ws.Range("Q5").Value = 65.45849407128499
$ws.Range("R5").Value = 589.1264466415649
$ws.Range("S5").Value = 0.004234760200889063
$ws.Range("T5").Value = 0.004234760200889064
# Row 6
$ws.Range("H6").Value = 6.766394999999999
$ws.Range("I6").Value = 0.01814126213173672
$ws.Range("J6").Value = 0.01814126213173672
$ws.Range("M6").Value = 92.64038833333332
$ws.Range("N6").Value = 277.921165
$ws.Range("O6").Value = 0.7451295270557885
$ws.Range("P6").Value = 0.7451295270557885
$ws.Range("Q6").Value = 208.9471534722416
$ws.Range("R6").Value = 1880.524381250175
$ws.Range("S6").Value = 0.01351759007241607
$ws.Range("T6").Value = 0.01351759007241607
# Row 7
$ws.Range("H7").Value = 6.766394999999999
$ws.Range("I7").Value = 0.01814126213173672
$ws.Range("J7").Value = 0.01814126213173672
$ws.Range("N7").Value = 6.766394999999999
$ws.Range("O7").Value = 0.01814126213173672
$ws.Range("P7").Value = 0.01814126213173672
$ws.Range("Q7").Value = 5.087122366224998
$ws.Range("R7").Value = 45.78410129602499
$ws.Range("S7").Value = 0.0003291053917323846
$ws.Range("T7").Value = 0.0003291053917323846
# Row 8
$ws.Range("H8").Value = 6.766394999999999
$ws.Range("I8").Value = 0.01814126213173672
$ws.Range("J8").Value = 0.01814126213173672
$ws.Range("M8").Value = 28.72545833333334
$ws.Range("N8").Value = 86.17637500000001
$ws.Range("O8").Value = 0.2310459570329316
$ws.Range("P8").Value = 0.2310459570329316
$ws.Range("Q8").Value = 64.78926587979166
$ws.Range("R8").Value = 583.103392918125
$ws.Range("S8").Value = 0.00419146527101239
$ws.Range("T8").Value = 0.004191465271012391
# Row 9
$ws.Range("H9").Value = 6.766394999999999
$ws.Range("I9").Value = 0.01814126213173672
$ws.Range("J9").Value = 0.01814126213173672
$ws.Range("M9").Value = 0.706587
$ws.Range("N9").Value = 2.119761
$ws.Range("O9").Value = 0.005683253779543222
$ws.Range("P9").Value = 0.005683253779543223
$ws.Range("Q9").Value = 1.593682247955
$ws.Range("R9").Value = 14.343140231595
$ws.Range("S9").Value = 0.000103101396575877
$ws.Range("T9").Value = 0.000103101396575877
# Row 10
$ws.Range("G10").Value = 28.72545833333334
$ws.Range("H10").Value = 86.17637500000001
$ws.Range("I10").Value = 0.2310459570329316
$ws.Range("J10").Value = 0.2310459570329316
$ws.Range("M10").Value = 92.64038833333332
$ws.Range("N10").Value = 277.921165
$ws.Range("O10").Value = 0.7451295270557885
$ws.Range("P10").Value = 0.7451295270557885
$ws.Range("Q10").Value = 2661.137615052986
$ws.Range("R10").Value = 23950.23853547687
$ws.Range("S10").Value = 0.1721591646921003
$ws.Range("T10").Value = 0.1721591646921004
# Row 11
$ws.Range("G11").Value = 28.72545833333334
$ws.Range("H11").Value = 86.17637500000001
$ws.Range("I11").Value = 0.2310459570329316
$ws.Range("J11").Value = 0.2310459570329316
$ws.Range("N11").Value = 6.766394999999999
$ws.Range("O11").Value = 0.01814126213173672
$ws.Range("P11").Value = 0.01814126213173672
$ws.Range("Q11").Value = 64.78926587979166
$ws.Range("R11").Value = 583.103392918125
$ws.Range("S11").Value = 0.00419146527101239
$ws.Range("T11").Value = 0.004191465271012391
# Row 12
$ws.Range("G12").Value = 28.72545833333334
$ws.Range("H12").Value = 86.17637500000001
$ws.Range("I12").Value = 0.2310459570329316
$ws.Range("J12").Value = 0.2310459570329316
$ws.Range("M12").Value = 28.72545833333334
$ws.Range("N12").Value = 86.17637500000001
$ws.Range("O12").Value = 0.2310459570329316
$ws.Range("P12").Value = 0.2310459570329316
$ws.Range("Q12").Value = 825.1519564600695
$ws.Range("R12").Value = 7426.367608140627
$ws.Range("S12").Value = 0.05338223426126326
$ws.Range("T12").Value = 0.05338223426126328
# Row 13
$ws.Range("G13").Value = 28.72545833333334
$ws.Range("H13").Value = 86.17637500000001
$ws.Range("I13").Value = 0.2310459570329316
$ws.Range("J13").Value = 0.2310459570329316
$ws.Range("M13").Value = 0.706587
$ws.Range("N13").Value = 2.119761
$ws.Range("O13").Value = 0.005683253779543222
$ws.Range("P13").Value = 0.005683253779543223
$ws.Range("Q13").Value = 20.297035427375
$ws.Range("R13").Value = 182.673318846375
$ws.Range("S13").Value = 0.001313092808555589
$ws.Range("T13").Value = 0.00131309280855559
# Row 14
$ws.Range("G14").Value = 0.706587
$ws.Range("H14").Value = 2.119761
$ws.Range("I14").Value = 0.005683253779543222
$ws.Range("J14").Value = 0.005683253779543223
$ws.Range("M14").Value = 92.64038833333332
$ws.Range("N14").Value = 277.921165
$ws.Range("O14").Value = 0.7451295270557885
$ws.Range("P14").Value = 0.7451295270557885
$ws.Range("Q14").Value = 65.45849407128499
$ws.Range("R14").Value = 589.1264466415649
$ws.Range("S14").Value = 0.004234760200889063
$ws.Range("T14").Value = 0.004234760200889064
# Row 15
$ws.Range("G15").Value = 0.706587
$ws.Range("H15").Value = 2.119761
$ws.Range("I15").Value = 0.005683253779543222
$ws.Range("J15").Value = 0.005683253779543223
$ws.Range("N15").Value = 6.766394999999999
$ws.Range("O15").Value = 0.01814126213173672
$ws.Range("P15").Value = 0.01814126213173672
$ws.Range("Q15").Value = 1.593682247955
$ws.Range("R15").Value = 14.343140231595
$ws.Range("S15").Value = 0.000103101396575877
$ws.Range("T15").Value = 0.000103101396575877
# Row 16
$ws.Range("G16").Value = 0.706587
$ws.Range("H16").Value = 2.119761
$ws.Range("I16").Value = 0.005683253779543222
$ws.Range("J16").Value = 0.005683253779543223
$ws.Range("M16").Value = 28.72545833333334
$ws.Range("N16").Value = 86.17637500000001
$ws.Range("O16").Value = 0.2310459570329316
$ws.Range("P16").Value = 0.2310459570329316
$ws.Range("Q16").Value = 20.297035427375
$ws.Range("R16").Value = 182.673318846375
$ws.Range("S16").Value = 0.001313092808555589
$ws.Range("T16").Value = 0.00131309280855559
# Row 17
$ws.Range("G17").Value = 0.706587
$ws.Range("H17").Value = 2.119761
$ws.Range("I17").Value = 0.005683253779543222
$ws.Range("J17").Value = 0.005683253779543223
$ws.Range("M17").Value = 0.706587
$ws.Range("N17").Value = 2.119761
$ws.Range("O17").Value = 0.005683253779543222
$ws.Range("P17").Value = 0.005683253779543223
$ws.Range("Q17").Value = 0.4992651885689999
$ws.Range("R17").Value = 4.493386697121
$ws.Range("S17").Value = 0.00003229937352269231
$ws.Range("T17").Value = 0.00003229937352269233
